$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unit change: kcal/mol -> kJ/mol (1 Hartree = 2625.5 kJ/mol, was 627.5095 kcal/mol) ---

# Header / label text updates (shared-string text edits)
$ws.Range("G1").Value  = "Relative Energy / kJ/mol"
$ws.Range("G9").Value  = "Relative to [B]_T0 / kJ/mol"
$ws.Range("H9").Value  = "Relative to [B-Trans]_T0 (gas phase) / kJ/mol"
$ws.Range("I9").Value  = "Relative to [B-Trans]_T0 (SMD) / kJ/mol"

# Formulas: swap the kcal/mol Hartree conversion factor for the kJ/mol one
$ws.Range("G2").Formula  = '=(D2-$D$3)*2625.5'
$ws.Range("G3").Formula  = '=(D3-$D$3)*2625.5'

$ws.Range("G5").Formula  = '=(D5-$D$8)*2625.5'
$ws.Range("G6").Formula  = '=(D6-$D$8)*2625.5'
$ws.Range("G7").Formula  = '=(D7-$D$8)*2625.5'
$ws.Range("G8").Formula  = '=(D8-$D$8)*2625.5'

$ws.Range("G10").Formula = '=((D10+D13+D12)-D2)*2625.5'
$ws.Range("H10").Formula = '=((E10+E12+E13)-E3)*2625.5'
$ws.Range("I10").Formula = '=((D10+D12+D13)-D3)*2625.5'

# --- View changes: zoom + selected cell ---
$ws.Activate() | Out-Null
$ws.Range("G9").Select() | Out-Null
$excel.ActiveWindow.Zoom = 114

# --- Column width changes (best achievable approximation of 22.1640625 / 10.83203125) ---
$ws.Columns.Item(1).ColumnWidth = 21.3
$ws.Columns.Item(7).ColumnWidth = 10.0
$ws.Columns.Item(8).ColumnWidth = 10.0
